# Update Name of Algo
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E7").Value = 12.1435
$ws.Range("B10").Value = 8.592200000000002
$ws.Range("B12").Value = 6.433599999999997
$ws.Range("D13").Value = -7.793800000000002
$ws.Range("B18").Value = 5.7458
$ws.Range("E20").Value = 13.0397
